# Update the "dSF" (F) column values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    4  = -1
    7  = 3
    9  = 5
    10 = 1
    12 = 0
    13 = -4
    15 = -7
    17 = -1
    21 = -1
    27 = 4
    41 = -3
    48 = 4
    52 = 0
    62 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
